$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and row 46/47 coin swap)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.333.22"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.933.83"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7492"
$ws.Range("E5").Value = "  +5.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.14"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.95"
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3184"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07027"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7818"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08042"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.928.45"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.405"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.12"
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.48"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.342.14"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.056"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.46"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000008000"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.179.13"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.704"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.536"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.70"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.08"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.226"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.531"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.422"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.142"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.343"
$ws.Range("E34").Value = "  +6.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05270"
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7578"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.785"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01960"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.808"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "79.02"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4508"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.985"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8367"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.755"
$ws.Range("E46").Value = "  +4.62%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.974"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.71"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.82"
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1236"
$ws.Range("E50").Value = "  +9.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "970.51"
$ws.Range("E51").Value = "  +6.55%  "

Write-Output "Applied cryptos update"
